$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the numeric/percent "price" and "volume" columns to remain text
# (they are stored as text strings like "314.15" / "2.13%", not numbers),
# matching the source data's inline-string representation.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Row 2-5: updated Price / Volume(1h) for existing coins ---
$ws.Range("D2").Value = "314.15"
$ws.Range("E2").Value = "2.13%"
$ws.Range("D3").Value = "40.82"
$ws.Range("E3").Value = "-0.55%"
$ws.Range("D4").Value = "5.160"
$ws.Range("E4").Value = "-1.53%"
$ws.Range("D5").Value = "0.07591"
$ws.Range("E5").Value = "-0.98%"

# --- Rows 6-18: the coin list shifted by one position (GateToken moved up
#     to row 6, each other coin shifted down one row), with refreshed
#     Price / Volume(1h) figures ---
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "4.327"
$ws.Range("E6").Value = "0.39%"

$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "1.679"
$ws.Range("E7").Value = "2.34%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9297"
$ws.Range("E8").Value = "1.56%"

$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "2.424"
$ws.Range("E9").Value = "-0.82%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1200"
$ws.Range("E10").Value = "-3.78%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1821"
$ws.Range("E11").Value = "-0.26%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09074"
$ws.Range("E12").Value = "-0.44%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.04139"
$ws.Range("E13").Value = "-2.97%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.1052"
$ws.Range("E14").Value = "0.12%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001281"
$ws.Range("E15").Value = "1.56%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005832"
$ws.Range("E16").Value = "1.29%"

$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "0.007522"
$ws.Range("E17").Value = "0.18%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.332"
$ws.Range("E18").Value = "-0.42%"

# --- Remaining rows: Price / Volume(1h) refreshes only ---
$ws.Range("D20").Value = "7.631"
$ws.Range("E20").Value = "4.31%"
$ws.Range("D23").Value = "0.04004"
$ws.Range("E23").Value = "-1.71%"
$ws.Range("D24").Value = "0.001283"
$ws.Range("E24").Value = "1.57%"
$ws.Range("D25").Value = "0.003973"
$ws.Range("E25").Value = "-7.37%"
$ws.Range("D26").Value = "0.0001352"
$ws.Range("E26").Value = "6.20%"
$ws.Range("D38").Value = "0.02410"
$ws.Range("E38").Value = "-2.49%"
$ws.Range("D39").Value = "0.05167"
$ws.Range("D40").Value = "0.007735"
$ws.Range("E40").Value = "-1.36%"
$ws.Range("D41").Value = "0.1298"
$ws.Range("E41").Value = "-1.10%"
$ws.Range("D42").Value = "0.007600"
$ws.Range("E42").Value = "10.43%"
$ws.Range("E43").Value = "72.51%"
$ws.Range("D44").Value = "0.008582"
$ws.Range("E44").Value = "12.47%"
$ws.Range("D45").Value = "0.3387"
$ws.Range("E45").Value = "10.62%"
$ws.Range("D46").Value = "0.00006597"
$ws.Range("E46").Value = "-1.96%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.15%"
$ws.Range("E48").Value = "-37.29%"
$ws.Range("E49").Value = "35.28%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.15%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.15%"

# Drop the temporary text number-format so the cells end up with the same
# (default/general) style they started with.
$ws.Range("D2:E51").ClearFormats()
